# Add a new "Account Funding Deposit" sheet after the existing sheets,
# populate it with the private-key / deposit-amount form, and make it
# the active tab.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the tab strip ---------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Account Funding Deposit"

# --- header row ----------------------------------------------------------
$ws.Range("A1").Value = "private key"
$ws.Range("B1").Value = "amount to deposit"

# --- data row --------------------------------------------------------------
$ws.Range("A2").Value = "fe1ef34ed4476ec0e7fabb2388d4a0e258d2ab28401a9836de60fcd44eb267b0"
$ws.Range("B2").Value = 1

# --- column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 69.16666666666667
$ws.Columns.Item(2).ColumnWidth = 18.022135416666668

# --- the new sheet becomes the active / selected tab ------------------------
$ws.Activate()
$ws.Range("B8").Select()
